$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(573,358,326,253,298,651,137,2648,1674,394,2140,390,3771,1705,516,979,1557,1084,369,256,168,611,983,353,317,1495,855,760,3201,1465,741,3745,1120,1788,643,147,1000,619,1067,329,2549,921,1488,105,769,233,729,5177,879,395,868,1850,1233,404,387,434,618,268,446,4275,448,318,1424,356,336,1817,1223,225,645,181,1468,436,1708,1714,396,2092,213,1381,329,304,726,851,1198,442,1107,248,809,1520,480,481,894,136,895,1951,2453,1123,212,908,551,711,545,821,873,112,388,498,262,332,481,798,2810,869,428,6194,753,294,240,683,1634,2435,572,1820,113,251,981,620,753,299,3040,314,2346,1054,867,379,1280,2090,677,279,394,351,488,915,1953,202,1772,889,880,1604,379,672,596,2135,166,974,554,259,409,331,228,761,437,439,488,1439,3392,1217,151,1616,290,508,406,7132,3965,1246,1932,375,1625,482,360,1985,183,973,1204,184,1578,868,380,1682,883,2515,474,1440,891,782,166,1897,108,149,260,1410,182,1420,6364,2131,1024,292,1958,819,412,2152,934,681,757,220,477,3197,1277,205,1266,3274,457,1943,1002,2502,399,402,755,2451,2155,583,380,2083,391,1511,887,543,332,291,418,723,6033,752,830,862,1539,836,1376,841,885,3709,166,529,3715,714,773,1120,1946,516,1510,1376,246,336,217,373,129,1096,9621,105,428,1166,153,637,189,584,989,606,3265,1611,2971,368,553,2316,1162,991,235,593,157,321,334,796,1327,756,6006,1408,110,2262,4205,2880,375,858,4670,2175,1303,765,5362,702,727,1667,314,2248,1129,135,207,273,651,4061,415,1593,1301,657,654,329)

$data = New-Object 'object[,]' $values.Length,1
for ($i = 0; $i -lt $values.Length; $i++) {
    $data[$i,0] = $values[$i]
}

$ws.Range("K2:K323").Value = $data

$ws.Range("K2:K323").Select()
